# Auto-generated edit script applying the cryptos.xlsx update
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '70.003.15'
$ws.Range('E2').Value = '  -1.38%  '
$ws.Range('D3').Value = '3.735.12'
$ws.Range('E3').Value = '  +1.51%  '
$ws.Range('E4').Value = '  +0.00%  '
$ws.Range('D5').Value = '''622.98'
$ws.Range('E5').Value = '  +0.37%  '
$ws.Range('D6').Value = '''180.66'
$ws.Range('E6').Value = '  -1.29%  '
$ws.Range('D7').Value = '3.733.12'
$ws.Range('E7').Value = '  +1.57%  '
$ws.Range('E8').Value = '  -0.01%  '
$ws.Range('E9').Value = '  -1.14%  '
$ws.Range('E10').Value = '  +1.98%  '
$ws.Range('D11').Value = '''6.31'
$ws.Range('E11').Value = '  -5.69%  '
$ws.Range('D12').Value = '''0.485'
$ws.Range('E12').Value = '  -3.52%  '
$ws.Range('D13').Value = '''41.00'
$ws.Range('E13').Value = '  +1.05%  '
$ws.Range('E14').Value = '  +1.60%  '
$ws.Range('D15').Value = '4.357.80'
$ws.Range('E15').Value = '  +1.51%  '
$ws.Range('D16').Value = '3.734.79'
$ws.Range('E16').Value = '  +1.58%  '
$ws.Range('D17').Value = '69.988.55'
$ws.Range('E17').Value = '  -1.46%  '
$ws.Range('E18').Value = '  -1.33%  '
$ws.Range('E19').Value = '  +0.03%  '
$ws.Range('D20').Value = '''16.80'
$ws.Range('E20').Value = '  -0.77%  '
$ws.Range('D21').Value = '''505.45'
$ws.Range('E21').Value = '  -2.79%  '
$ws.Range('D22').Value = '''9.41'
$ws.Range('E22').Value = '  +1.41%  '
$ws.Range('D23').Value = '''0.722'
$ws.Range('E23').Value = '  -2.95%  '
$ws.Range('D24').Value = '''2.52'
$ws.Range('E24').Value = '  -0.82%  '
$ws.Range('D25').Value = '''86.63'
$ws.Range('E25').Value = '  -2.27%  '
$ws.Range('E26').Value = '  -3.24%  '
$ws.Range('E27').Value = '  +1.21%  '
$ws.Range('E28').Value = '  +22.79%  '
$ws.Range('E29').Value = '  -0.14%  '
$ws.Range('E30').Value = '  -2.58%  '
$ws.Range('E31').Value = '  +0.17%  '
$ws.Range('D32').Value = '''7.94'
$ws.Range('E32').Value = '  -3.19%  '
$ws.Range('D33').Value = '''31.15'
$ws.Range('E33').Value = '  -1.70%  '
$ws.Range('E34').Value = '  -0.55%  '
$ws.Range('E35').Value = '  -0.02%  '
$ws.Range('E36').Value = '  +2.39%  '
$ws.Range('D37').Value = '''6.20'
$ws.Range('B38').Value = 'Kaspa'
$ws.Range('C38').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D38').Value = '''0.135'
$ws.Range('E38').Value = '  +2.40%  '
$ws.Range('B39').Value = 'TheGraph'
$ws.Range('C39').Value = 'https://coinranking.com/coin/qhd1biQ7M+thegraph-grt'
$ws.Range('D39').Value = '''0.335'
$ws.Range('E39').Value = '  -3.74%  '
$ws.Range('E40').Value = '  -4.89%  '
$ws.Range('D41').Value = '''50.43'
$ws.Range('E41').Value = '  -2.03%  '
$ws.Range('D42').Value = '''45.00'
$ws.Range('E42').Value = '  -1.50%  '
$ws.Range('D43').Value = '''425.98'
$ws.Range('E43').Value = '  -2.13%  '
$ws.Range('D44').Value = '''8.71'
$ws.Range('E44').Value = '  -1.43%  '
$ws.Range('E45').Value = '  -1.03%  '
$ws.Range('D46').Value = '2.998.00'
$ws.Range('E46').Value = '  -3.88%  '
$ws.Range('D47').Value = '''0.0363'
$ws.Range('E47').Value = '  -1.82%  '
$ws.Range('D48').Value = '''27.32'
$ws.Range('E48').Value = '  -3.66%  '
$ws.Range('E49').Value = '  -0.05%  '
$ws.Range('D50').Value = '''137.39'
$ws.Range('E50').Value = '  -1.94%  '
$ws.Range('E51').Value = '  +1.56%  '
